$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "66.172.93"
$ws.Range("E2").Value2 = "  -1.71%  "
$ws.Range("D3").Value2 = "2.536.61"
$ws.Range("E3").Value2 = "  -3.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.00"
$ws.Range("E4").Value2 = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "579.58"
$ws.Range("E5").Value2 = "  -2.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "168.09"
$ws.Range("E6").Value2 = "  +0.04%  "
$ws.Range("E7").Value2 = "  +0.01%  "
$ws.Range("E8").Value2 = "  -1.78%  "
$ws.Range("D9").Value2 = "2.535.59"
$ws.Range("E9").Value2 = "  -3.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.139"
$ws.Range("E10").Value2 = "  -0.02%  "
$ws.Range("E11").Value2 = "  -0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.350"
$ws.Range("E12").Value2 = "  -3.73%  "
$ws.Range("E13").Value2 = "  -2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "26.49"
$ws.Range("E14").Value2 = "  -4.74%  "
$ws.Range("D15").Value2 = "3.001.28"
$ws.Range("E15").Value2 = "  -3.30%  "
$ws.Range("E16").Value2 = "  -3.33%  "
$ws.Range("D17").Value2 = "66.123.31"
$ws.Range("E17").Value2 = "  -1.81%  "
$ws.Range("D18").Value2 = "2.515.52"
$ws.Range("E18").Value2 = "  -4.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "11.31"
$ws.Range("E19").Value2 = "  -6.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "7.62"
$ws.Range("E20").Value2 = "  -5.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "346.78"
$ws.Range("E21").Value2 = "  -2.83%  "
$ws.Range("E22").Value2 = "  -3.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "4.55"
$ws.Range("E23").Value2 = "  -2.99%  "
$ws.Range("E24").Value2 = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "1.92"
$ws.Range("E25").Value2 = "  -0.63%  "
$ws.Range("E26").Value2 = "  -0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "9.89"
$ws.Range("E27").Value2 = "  -5.45%  "
$ws.Range("E28").Value2 = "  -2.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "1.00"
$ws.Range("E29").Value2 = "  -0.12%  "
$ws.Range("E30").Value2 = "  -3.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "524.39"
$ws.Range("E31").Value2 = "  -4.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "8.16"
$ws.Range("E32").Value2 = "  +2.74%  "
$ws.Range("E33").Value2 = "  -3.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.81"
$ws.Range("E34").Value2 = "  -5.19%  "
$ws.Range("E35").Value2 = "  -3.48%  "
$ws.Range("E36").Value2 = "  -0.09%  "
$ws.Range("B37").Value2 = "Monero"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "156.82"
$ws.Range("E37").Value2 = "  -0.98%  "
$ws.Range("B38").Value2 = "ImmutableX"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "1.44"
$ws.Range("E38").Value2 = "  -4.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "18.67"
$ws.Range("E39").Value2 = "  -1.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "18.28"
$ws.Range("E40").Value2 = "  +0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.354"
$ws.Range("E41").Value2 = "  -3.38%  "
$ws.Range("E42").Value2 = "  -2.73%  "
$ws.Range("E43").Value2 = "  -2.83%  "
$ws.Range("E44").Value2 = "  +0.03%  "
$ws.Range("E45").Value2 = "  -0.87%  "
$ws.Range("B46").Value2 = "Aave"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "147.33"
$ws.Range("E46").Value2 = "  -3.16%  "
$ws.Range("B47").Value2 = "BabyDogeCoin"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value2 = "0.0₆0282"
$ws.Range("E47").Value2 = "  -4.63%  "
$ws.Range("E48").Value2 = "  -4.15%  "
$ws.Range("E49").Value2 = "  -2.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.71"
$ws.Range("E50").Value2 = "  -0.11%  "
$ws.Range("E51").Value2 = "  -1.85%  "
